# lesson_25_Food ko - part 2
#
# The original single run
#     "I have a c……………………..for roast meet"
# is corrected to "...roast meat" and split into three runs:
#     "I "  /  "have a c……………………..for roast mea"  /  "t"
# with the document's "_GoBack" bookmark now sitting between "mea" and "t"
# (it used to sit between "...tough mea" and "t" a bit further down the
# document - Word only ever keeps one _GoBack, so it is removed from its
# old spot and re-added at the new one).

$d = $word.ActiveDocument

# --- Step 0: drop the existing _GoBack bookmark (currently inside the
#     "tough mea|t" run pair) -------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 1: find the sentence and fix the "meet" -> "meat" typo -------
$sentence = $d.Content
$found = $sentence.Find.Execute(
    "I have a c……………………..for roast meet",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the target sentence 'I have a c……………………..for roast meet'"
}

$sentenceStart = $sentence.Start
$sentenceEnd   = $sentence.End

# the second "e" of "meet" is the character right before the final "t"
$typo = $d.Range($sentenceEnd - 2, $sentenceEnd - 1)
$typo.Text = "a"

# --- Step 2: split "I " off into its own run ----------------------------
# Inserting (and immediately removing) a bookmark at a position forces the
# run to split there without altering any text.
$splitPoint = $d.Range($sentenceStart + 2, $sentenceStart + 2)
$d.Bookmarks.Add("TempRunSplit", $splitPoint)
$d.Bookmarks("TempRunSplit").Delete()

# --- Step 3: re-add _GoBack between "mea" and "t" of "roast meat" ------
$newMarkPos = $sentenceEnd - 1
$newMark = $d.Range($newMarkPos, $newMarkPos)
$d.Bookmarks.Add("_GoBack", $newMark)
